$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Featured Products" translation was retranslated from "המיוחדים" to "ספיישלים"
$ws.Range("C124").Value = "ספיישלים"

# Two new key/value translation rows appended at the bottom of the table
$ws.Range("B149").Value = "Self picking"
$ws.Range("C149").Value = "עצמי"
$ws.Range("B150").Value = "Delivery"
$ws.Range("C150").Value = "משלוח"

# Match the formatting used by the most-recently-added rows above (green
# JetBrains Mono font, style index 3) instead of leaving default formatting.
$ws.Range("B133").Copy()
$ws.Range("B149").PasteSpecial(-4122)
$ws.Range("B150").PasteSpecial(-4122)

# Update the active selection to reflect where editing left off.
$ws.Range("C154").Select()
